# Updates crypto price/volume data per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'41.421.71"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').Value = "'2.215.25"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'248.94"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.31%  '
$ws.Range('D6').Value = "'0.624"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.25%  '
$ws.Range('D7').Value = "'70.12"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.24%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.568"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.93%  '
$ws.Range('D10').Value = "'41.21"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +15.69%  '
$ws.Range('D11').Value = "'0.0955"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.83%  '
$ws.Range('D12').Value = "'58.59"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = "'6.97"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.39%  '
$ws.Range('D15').Value = "'2.549.85"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').Value = "'14.76"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('D17').Value = "'0.848"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('D18').Value = "'2.214.57"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').Value = "'41.442.31"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('D20').Value = "'0.0₃0959"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.67%  '
$ws.Range('D21').Value = "'6.17"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.85%  '
$ws.Range('D22').Value = "'72.31"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.12%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = "'232.87"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = "'2.21"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.29%  '
$ws.Range('D25').Value = "'3.85"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.90%  '
$ws.Range('D27').Value = "'2.48"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.78%  '
$ws.Range('D28').Value = "'10.45"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.68%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Value = "'170.62"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.30%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = "'2.10"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.62%  '
$ws.Range('D31').Value = "'20.49"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').Value = "'0.119"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('E33').Value = '  -2.61%  '
$ws.Range('D34').Value = "'5.49"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.17%  '
$ws.Range('D35').Value = "'0.0714"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').Value = "'4.64"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').Value = "'25.83"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +17.44%  '
$ws.Range('D38').Value = "'3.93"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.62%  '
$ws.Range('D39').Value = "'0.0288"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.81%  '
$ws.Range('D40').Value = "'2.27"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('D41').Value = "'68.23"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.44%  '
$ws.Range('D42').Value = "'5.89"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('D43').Value = "'11.80"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +19.72%  '
$ws.Range('D44').Value = "'0.206"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.47%  '
$ws.Range('E45').Value = '  -1.69%  '
$ws.Range('D46').Value = "'8.69"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.03%  '
$ws.Range('D47').Value = "'4.73"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.24%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').Value = "'1.15"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.22%  '
$ws.Range('D51').Value = "'1.18"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.36%  '
